$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "638×5=3190"; New = "834×3=2502" },
    @{ Old = "338×9=3042"; New = "248×9=2232" },
    @{ Old = "380×9=3420"; New = "689×8=5512" },
    @{ Old = "693×9=6237"; New = "131×8=1048" },
    @{ Old = "914×4=3656"; New = "785×7=5495" },
    @{ Old = "629×8=5032"; New = "561×8=4488" },
    @{ Old = "976×9=8784"; New = "194×4=776" },
    @{ Old = "117×7=819";  New = "706×3=2118" },
    @{ Old = "387×6=2322"; New = "751×4=3004" },
    @{ Old = "805×8=6440"; New = "548×7=3836" },
    @{ Old = "771×4=3084"; New = "956×6=5736" },
    @{ Old = "136×7=952";  New = "968×2=1936" },
    @{ Old = "507×7=3549"; New = "680×4=2720" },
    @{ Old = "110×2=220";  New = "807×7=5649" },
    @{ Old = "730×9=6570"; New = "509×4=2036" },
    @{ Old = "241×8=1928"; New = "222×5=1110" },
    @{ Old = "248×7=1736"; New = "803×6=4818" },
    @{ Old = "820×4=3280"; New = "318×3=954" },
    @{ Old = "835×8=6680"; New = "312×8=2496" },
    @{ Old = "384×2=768";  New = "653×2=1306" },
    @{ Old = "958×8=7664"; New = "357×6=2142" },
    @{ Old = "261×2=522";  New = "838×3=2514" },
    @{ Old = "329×7=2303"; New = "544×2=1088" },
    @{ Old = "523×5=2615"; New = "954×6=5724" },
    @{ Old = "855×3=2565"; New = "679×4=2716" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
